$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 171; this shifts existing rows 171-204 down to 172-205.
$ws.Rows("171:171").Insert()

# Populate the new row 171 with the weekly price entry.
$ws.Range("A171").Value = 8
$ws.Range("B171").Value = "Terminal La Palmera de La Serena"
$ws.Range("C171").Value = "Coquimbo"
$ws.Range("D171").NumberFormat = $ws.Range("D172").NumberFormat
$ws.Range("D171").Value = 45258
$ws.Range("E171").Value = 4
$ws.Range("F171").Value = 100114007
$ws.Range("G171").Value = "Jengibre"
$ws.Range("H171").Value = "Sin especificar"
$ws.Range("I171").Value = "Primera"
$ws.Range("J171").Value = 360
$ws.Range("K171").Value = 23000
$ws.Range("L171").Value = 24000
$ws.Range("M171").Value = 23500
$ws.Range("N171").Value = "`$/caja 13 kilos"
$ws.Range("O171").Value = "Perú"
$ws.Range("P171").Value = 1808
$ws.Range("Q171").Value = 13
$ws.Range("R171").Value = "Hortaliza"
